# Feat: Added both UIs
# Populates the "Translation" sheet with the single-use text id rows, and
# sets the Widget Wildcard Characters cell (H4) on the "Typography" sheet.

$wb = $excel.ActiveWorkbook
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# Typography sheet: Widget Wildcard Characters for the first (Default) font entry
$wsTypography.Range("H4").Value = "0123456789 :APM"

# Translation sheet: header already present in row 3 (TEXT ID / TYPOGRAPHY NAME / ALIGNMENT / DIRECTION / GB)
# Row 4 - first entry
$wsTranslation.Range("B4").Value = "SingleUseId1"
$wsTranslation.Range("C4").Value = "Default"
$wsTranslation.Range("D4").Value = "Center"
$wsTranslation.Range("E4").Value = "LTR"
$wsTranslation.Range("F4").Value = ""
$wsTranslation.Range("F4").Value = "<>"

# Row 5 - was drafted then edited afterwards
$wsTranslation.Range("B5").Value = "SingleUseId2"
$wsTranslation.Range("C5").Value = "Default"
$wsTranslation.Range("D5").Value = "Right"
$wsTranslation.Range("E5").Value = "LTR"
$wsTranslation.Range("F5").Value = "0"
$wsTranslation.Range("B5").Value = "SingleUseId3"
$wsTranslation.Range("F5").Value = "<value>"

# Row 6
$wsTranslation.Range("B6").Value = "SingleUseId4"
$wsTranslation.Range("C6").Value = "Default"
$wsTranslation.Range("D6").Value = "Center"
$wsTranslation.Range("E6").Value = "LTR"
$wsTranslation.Range("F6").Value = "1"

# Row 7
$wsTranslation.Range("B7").Value = "SingleUseId5"
$wsTranslation.Range("C7").Value = "Default"
$wsTranslation.Range("D7").Value = "Center"
$wsTranslation.Range("E7").Value = "LTR"
$wsTranslation.Range("F7").Value = "2"

# Row 8
$wsTranslation.Range("B8").Value = "SingleUseId6"
$wsTranslation.Range("C8").Value = "Default"
$wsTranslation.Range("D8").Value = "Center"
$wsTranslation.Range("E8").Value = "LTR"
$wsTranslation.Range("F8").Value = "3"

# Row 9
$wsTranslation.Range("B9").Value = "SingleUseId7"
$wsTranslation.Range("C9").Value = "Default"
$wsTranslation.Range("D9").Value = "Center"
$wsTranslation.Range("E9").Value = "LTR"
$wsTranslation.Range("F9").Value = "0"

# Row 10
$wsTranslation.Range("B10").Value = "SingleUseId8"
$wsTranslation.Range("C10").Value = "Default"
$wsTranslation.Range("D10").Value = "Center"
$wsTranslation.Range("E10").Value = "LTR"
$wsTranslation.Range("F10").Value = "4"

# Row 11
$wsTranslation.Range("B11").Value = "SingleUseId9"
$wsTranslation.Range("C11").Value = "Default"
$wsTranslation.Range("D11").Value = "Center"
$wsTranslation.Range("E11").Value = "LTR"
$wsTranslation.Range("F11").Value = "5"

# Row 12
$wsTranslation.Range("B12").Value = "SingleUseId11"
$wsTranslation.Range("C12").Value = "Default"
$wsTranslation.Range("D12").Value = "Center"
$wsTranslation.Range("E12").Value = "LTR"
$wsTranslation.Range("F12").Value = "6"

# Row 13
$wsTranslation.Range("B13").Value = "SingleUseId12"
$wsTranslation.Range("C13").Value = "Default"
$wsTranslation.Range("D13").Value = "Center"
$wsTranslation.Range("E13").Value = "LTR"
$wsTranslation.Range("F13").Value = "C"

# Row 14
$wsTranslation.Range("B14").Value = "SingleUseId13"
$wsTranslation.Range("C14").Value = "Default"
$wsTranslation.Range("D14").Value = "Center"
$wsTranslation.Range("E14").Value = "LTR"
$wsTranslation.Range("F14").Value = "7"

# Row 15
$wsTranslation.Range("B15").Value = "SingleUseId14"
$wsTranslation.Range("C15").Value = "Default"
$wsTranslation.Range("D15").Value = "Center"
$wsTranslation.Range("E15").Value = "LTR"
$wsTranslation.Range("F15").Value = "8"

# Row 16
$wsTranslation.Range("B16").Value = "SingleUseId15"
$wsTranslation.Range("C16").Value = "Default"
$wsTranslation.Range("D16").Value = "Center"
$wsTranslation.Range("E16").Value = "LTR"
$wsTranslation.Range("F16").Value = "9"

# Row 17
$wsTranslation.Range("B17").Value = "SingleUseId16"
$wsTranslation.Range("C17").Value = "Default"
$wsTranslation.Range("D17").Value = "Center"
$wsTranslation.Range("E17").Value = "LTR"
$wsTranslation.Range("F17").Value = "E"

# Revisit row 5's F cell with its final value
$wsTranslation.Range("F5").Value = "Bomb<value>"

# Row 18
$wsTranslation.Range("B18").Value = "SingleUseId17"
$wsTranslation.Range("C18").Value = "Default"
$wsTranslation.Range("D18").Value = "Center"
$wsTranslation.Range("E18").Value = "LTR"
$wsTranslation.Range("F18").Value = "Bomb<value>"

# Row 19
$wsTranslation.Range("B19").Value = "SingleUseId18"
$wsTranslation.Range("C19").Value = "Default"
$wsTranslation.Range("D19").Value = "Center"
$wsTranslation.Range("E19").Value = "LTR"
$wsTranslation.Range("F19").Value = "End<value>"
